$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 12.73908466666667
$ws.Range("H2").Value = 38.217254
$ws.Range("I2").Value = 0.005953388968763418
$ws.Range("J2").Value = 0.006105597140986208
$ws.Range("M2").Value = 121.928739
$ws.Range("N2").Value = 365.786217
$ws.Range("O2").Value = 0.2282232151508951
$ws.Range("P2").Value = 0.2419720431319445
$ws.Range("Q2").Value = 1553.260529420902
$ws.Range("R2").Value = 13979.34476478812
$ws.Range("S2").Value = 0.001358701571495059
$ws.Range("T2").Value = 0.001477383814744992
$ws.Range("G3").Value = 12.73908466666667
$ws.Range("H3").Value = 38.217254
$ws.Range("I3").Value = 0.005953388968763418
$ws.Range("J3").Value = 0.006105597140986208
$ws.Range("O3").Value = 0.2768624053389947
$ws.Range("P3").Value = 0.2935413991166814
$ws.Range("Q3").Value = 1884.293173283273
$ws.Range("R3").Value = 16958.63855954946
$ws.Range("S3").Value = 0.001648269589810477
$ws.Range("T3").Value = 0.001792245527207902
$ws.Range("G4").Value = 12.73908466666667
$ws.Range("H4").Value = 38.217254
$ws.Range("I4").Value = 0.005953388968763418
$ws.Range("J4").Value = 0.006105597140986208
$ws.Range("M4").Value = 83.50496933333334
$ws.Range("N4").Value = 250.514908
$ws.Range("O4").Value = 0.1563025480180701
$ws.Range("P4").Value = 0.1657186665504434
$ws.Range("Q4").Value = 1063.776874424737
$ws.Range("R4").Value = 9573.991869822632
$ws.Range("S4").Value = 0.0009305298651603927
$ws.Range("T4").Value = 0.001011811416698434
$ws.Range("G5").Value = 12.73908466666667
$ws.Range("H5").Value = 38.217254
$ws.Range("I5").Value = 0.005953388968763418
$ws.Range("J5").Value = 0.006105597140986208
$ws.Range("M5").Value = 91.06846250000001
$ws.Range("N5").Value = 182.136925
$ws.Range("O5").Value = 0.1704597085236707
$ws.Range("P5").Value = 0.1204857969594293
$ws.Range("Q5").Value = 1160.128854250659
$ws.Range("R5").Value = 6960.773125503952
$ws.Range("S5").Value = 0.001014812948343449
$ws.Range("T5").Value = 0.0007356377374449365
$ws.Range("G6").Value = 12.73908466666667
$ws.Range("H6").Value = 38.217254
$ws.Range("I6").Value = 0.005953388968763418
$ws.Range("J6").Value = 0.006105597140986208
$ws.Range("M6").Value = 89.83562999999999
$ws.Range("N6").Value = 269.50689
$ws.Range("O6").Value = 0.1681521229683693
$ws.Range("P6").Value = 0.1782820942415013
$ws.Range("Q6").Value = 1144.42369665334
$ws.Range("R6").Value = 10299.81326988006
$ws.Range("S6").Value = 0.00100107499395404
$ws.Range("T6").Value = 0.001088518644889944
$ws.Range("I7").Value = 0.9182810852447438
$ws.Range("J7").Value = 0.9417584502053091
$ws.Range("M7").Value = 121.928739
$ws.Range("N7").Value = 365.786217
$ws.Range("O7").Value = 0.2282232151508951
$ws.Range("P7").Value = 0.2419720431319445
$ws.Range("Q7").Value = 239582.8278831099
$ws.Range("R7").Value = 2156245.450947989
$ws.Range("S7").Value = 0.2095730616868086
$ws.Range("T7").Value = 0.2278792163329523
$ws.Range("I8").Value = 0.9182810852447438
$ws.Range("J8").Value = 0.9417584502053091
$ws.Range("O8").Value = 0.2768624053389947
$ws.Range("P8").Value = 0.2935413991166814
$ws.Range("Q8").Value = 290642.9916070527
$ws.Range("S8").Value = 0.2542375100381622
$ws.Range("T8").Value = 0.276445093103224
$ws.Range("I9").Value = 0.9182810852447438
$ws.Range("J9").Value = 0.9417584502053091
$ws.Range("M9").Value = 83.50496933333334
$ws.Range("N9").Value = 250.514908
$ws.Range("O9").Value = 0.1563025480180701
$ws.Range("P9").Value = 0.1657186665504434
$ws.Range("Q9").Value = 164082.3718776619
$ws.Range("R9").Value = 1476741.346898957
$ws.Range("S9").Value = 0.1435296734205521
$ws.Range("T9").Value = 0.156066954580636
$ws.Range("I10").Value = 0.9182810852447438
$ws.Range("J10").Value = 0.9417584502053091
$ws.Range("M10").Value = 91.06846250000001
$ws.Range("N10").Value = 182.136925
$ws.Range("O10").Value = 0.1704597085236707
$ws.Range("P10").Value = 0.1204857969594293
$ws.Range("Q10").Value = 178944.1927773645
$ws.Range("R10").Value = 1073665.156664187
$ws.Range("S10").Value = 0.1565299261336191
$ws.Range("T10").Value = 0.1134685174162637
$ws.Range("I11").Value = 0.9182810852447438
$ws.Range("J11").Value = 0.9417584502053091
$ws.Range("M11").Value = 89.83562999999999
$ws.Range("N11").Value = 269.50689
$ws.Range("O11").Value = 0.1681521229683693
$ws.Range("P11").Value = 0.1782820942415013
$ws.Range("Q11").Value = 176521.749151041
$ws.Range("R11").Value = 1588695.742359369
$ws.Range("S11").Value = 0.1544109139656018
$ws.Range("T11").Value = 0.1678986687722331
$ws.Range("G12").Value = 1.091866333333334
$ws.Range("H12").Value = 3.275599000000001
$ws.Range("I12").Value = 0.0005102646818291153
$ws.Range("J12").Value = 0.0005233104369407934
$ws.Range("M12").Value = 121.928739
$ws.Range("N12").Value = 365.786217
$ws.Range("O12").Value = 0.2282232151508951
$ws.Range("P12").Value = 0.2419720431319445
$ws.Range("Q12").Value = 133.129885179887
$ws.Range("R12").Value = 1198.168966618983
$ws.Range("S12").Value = 0.0001164542462649892
$ws.Range("T12").Value = 0.0001266264956188344
$ws.Range("G13").Value = 1.091866333333334
$ws.Range("H13").Value = 3.275599000000001
$ws.Range("I13").Value = 0.0005102646818291153
$ws.Range("J13").Value = 0.0005233104369407934
$ws.Range("O13").Value = 0.2768624053389947
$ws.Range("P13").Value = 0.2935413991166814
$ws.Range("Q13").Value = 161.5026771445567
$ws.Range("R13").Value = 1453.52409430101
$ws.Range("S13").Value = 0.0001412731071707457
$ws.Range("T13").Value = 0.0001536132778319624
$ws.Range("G14").Value = 1.091866333333334
$ws.Range("H14").Value = 3.275599000000001
$ws.Range("I14").Value = 0.0005102646818291153
$ws.Range("J14").Value = 0.0005233104369407934
$ws.Range("M14").Value = 83.50496933333334
$ws.Range("N14").Value = 250.514908
$ws.Range("O14").Value = 0.1563025480180701
$ws.Range("P14").Value = 0.1657186665504434
$ws.Range("Q14").Value = 91.17626468109914
$ws.Range("R14").Value = 820.5863821298922
$ws.Range("S14").Value = 0.00007975566993352054
$ws.Range("T14").Value = 0.00008672230780175816
$ws.Range("G15").Value = 1.091866333333334
$ws.Range("H15").Value = 3.275599000000001
$ws.Range("I15").Value = 0.0005102646818291153
$ws.Range("J15").Value = 0.0005233104369407934
$ws.Range("M15").Value = 91.06846250000001
$ws.Range("N15").Value = 182.136925
$ws.Range("O15").Value = 0.1704597085236707
$ws.Range("P15").Value = 0.1204857969594293
$ws.Range("Q15").Value = 99.43458823217921
$ws.Range("R15").Value = 596.6075293930752
$ws.Range("S15").Value = 0.00008697956893451458
$ws.Range("T15").Value = 0.00006305147505199868
$ws.Range("G16").Value = 1.091866333333334
$ws.Range("H16").Value = 3.275599000000001
$ws.Range("I16").Value = 0.0005102646818291153
$ws.Range("J16").Value = 0.0005233104369407934
$ws.Range("M16").Value = 89.83562999999999
$ws.Range("N16").Value = 269.50689
$ws.Range("O16").Value = 0.1681521229683693
$ws.Range("P16").Value = 0.1782820942415013
$ws.Range("Q16").Value = 98.08849993079002
$ws.Range("R16").Value = 882.7964993771102
$ws.Range("S16").Value = 0.00008580208952534522
$ws.Range("T16").Value = 0.00009329688063623975
$ws.Range("G17").Value = 160.0313415
$ws.Range("H17").Value = 320.062683
$ws.Range("I17").Value = 0.0747878554913321
$ws.Range("J17").Value = 0.05113328661083746
$ws.Range("M17").Value = 121.928739
$ws.Range("N17").Value = 365.786217
$ws.Range("O17").Value = 0.2282232151508951
$ws.Range("P17").Value = 0.2419720431319445
$ws.Range("Q17").Value = 19512.41966957337
$ws.Range("R17").Value = 117074.5180174402
$ws.Range("S17").Value = 0.01706832483447234
$ws.Range("T17").Value = 0.01237282583327565
$ws.Range("G18").Value = 160.0313415
$ws.Range("H18").Value = 320.062683
$ws.Range("I18").Value = 0.0747878554913321
$ws.Range("J18").Value = 0.05113328661083746
$ws.Range("O18").Value = 0.2768624053389947
$ws.Range("P18").Value = 0.2935413991166814
$ws.Range("Q18").Value = 23670.92865697369
$ws.Range("R18").Value = 142025.5719418422
$ws.Range("S18").Value = 0.02070594556147535
$ws.Range("T18").Value = 0.0150097364931795
$ws.Range("G19").Value = 160.0313415
$ws.Range("H19").Value = 320.062683
$ws.Range("I19").Value = 0.0747878554913321
$ws.Range("J19").Value = 0.05113328661083746
$ws.Range("M19").Value = 83.50496933333334
$ws.Range("N19").Value = 250.514908
$ws.Range("O19").Value = 0.1563025480180701
$ws.Range("P19").Value = 0.1657186665504434
$ws.Range("Q19").Value = 13363.41226432969
$ws.Range("R19").Value = 80180.47358597816
$ws.Range("S19").Value = 0.01168953237410242
$ws.Range("T19").Value = 0.008473740073489625
$ws.Range("G20").Value = 160.0313415
$ws.Range("H20").Value = 320.062683
$ws.Range("I20").Value = 0.0747878554913321
$ws.Range("J20").Value = 0.05113328661083746
$ws.Range("M20").Value = 91.06846250000001
$ws.Range("N20").Value = 182.136925
$ws.Range("O20").Value = 0.1704597085236707
$ws.Range("P20").Value = 0.1204857969594293
$ws.Range("Q20").Value = 14573.80822221745
$ws.Range("R20").Value = 58295.23288886978
$ws.Range("S20").Value = 0.01274831604816288
$ws.Range("T20").Value = 0.006160834788461669
$ws.Range("G21").Value = 160.0313415
$ws.Range("H21").Value = 320.062683
$ws.Range("I21").Value = 0.0747878554913321
$ws.Range("J21").Value = 0.05113328661083746
$ws.Range("M21").Value = 89.83562999999999
$ws.Range("N21").Value = 269.50689
$ws.Range("O21").Value = 0.1681521229683693
$ws.Range("P21").Value = 0.1782820942415013
$ws.Range("Q21").Value = 14376.51638339764
$ws.Range("R21").Value = 86259.09830038587
$ws.Range("S21").Value = 0.01257573667311911
$ws.Range("T21").Value = 0.009116149422431021
$ws.Range("G22").Value = 1.000156333333333
$ws.Range("H22").Value = 3.000469
$ws.Range("I22").Value = 0.0004674056133315229
$ws.Range("J22").Value = 0.0004793556059265206
$ws.Range("M22").Value = 121.928739
$ws.Range("N22").Value = 365.786217
$ws.Range("O22").Value = 0.2282232151508951
$ws.Range("P22").Value = 0.2419720431319445
$ws.Range("Q22").Value = 121.947800526197
$ws.Range("R22").Value = 1097.530204735773
$ws.Range("S22").Value = 0.0001066728118540962
$ws.Range("T22").Value = 0.0001159906553527914
$ws.Range("G23").Value = 1.000156333333333
$ws.Range("H23").Value = 3.000469
$ws.Range("I23").Value = 0.0004674056133315229
$ws.Range("J23").Value = 0.0004793556059265206
$ws.Range("O23").Value = 0.2768624053389947
$ws.Range("P23").Value = 0.2935413991166814
$ws.Range("Q23").Value = 147.9374539402566
$ws.Range("R23").Value = 1331.43708546231
$ws.Range("S23").Value = 0.0001294070423759135
$ws.Range("T23").Value = 0.0001407107152380954
$ws.Range("G24").Value = 1.000156333333333
$ws.Range("H24").Value = 3.000469
$ws.Range("I24").Value = 0.0004674056133315229
$ws.Range("J24").Value = 0.0004793556059265206
$ws.Range("M24").Value = 83.50496933333334
$ws.Range("N24").Value = 250.514908
$ws.Range("O24").Value = 0.1563025480180701
$ws.Range("P24").Value = 0.1657186665504434
$ws.Range("Q24").Value = 83.51802394353911
$ws.Range("R24").Value = 751.662215491852
$ws.Range("S24").Value = 0.00007305668832166585
$ws.Range("T24").Value = 0.0000794381718176228
$ws.Range("G25").Value = 1.000156333333333
$ws.Range("H25").Value = 3.000469
$ws.Range("I25").Value = 0.0004674056133315229
$ws.Range("J25").Value = 0.0004793556059265206
$ws.Range("M25").Value = 91.06846250000001
$ws.Range("N25").Value = 182.136925
$ws.Range("O25").Value = 0.1704597085236707
$ws.Range("P25").Value = 0.1204857969594293
$ws.Range("Q25").Value = 91.08269953630416
$ws.Range("R25").Value = 546.4961972178251
$ws.Range("S25").Value = 0.00007967382461081895
$ws.Range("T25").Value = 0.00005775554220702698
$ws.Range("G26").Value = 1.000156333333333
$ws.Range("H26").Value = 3.000469
$ws.Range("I26").Value = 0.0004674056133315229
$ws.Range("J26").Value = 0.0004793556059265206
$ws.Range("M26").Value = 89.83562999999999
$ws.Range("N26").Value = 269.50689
$ws.Range("O26").Value = 0.1681521229683693
$ws.Range("P26").Value = 0.1782820942415013
$ws.Range("Q26").Value = 89.84967430348998
$ws.Range("R26").Value = 808.6470687314099
$ws.Range("S26").Value = 0.00007859524616902832
$ws.Range("T26").Value = 0.00008546052131098389
